$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gens")
$ws.Range("C8").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("C25").Value = 400

$ws = $wb.Worksheets.Item("lines")
$ws.Range("C2").Value = 24.990052
$ws.Range("D2").Value = 0.67236251
$ws.Range("C3").Value = -91.732021
$ws.Range("D3").Value = -2.555392
$ws.Range("C4").Value = -41.25803
$ws.Range("D4").Value = -2.5430587
$ws.Range("C5").Value = -50.509948
$ws.Range("D5").Value = -2.1918025
$ws.Range("D6").Value = 12.534563
$ws.Range("C7").Value = 2.8699097
$ws.Range("D7").Value = -0.75890982
$ws.Range("C8").Value = -274.60193
$ws.Range("D8").Value = -0.48161142
$ws.Range("C9").Value = -124.50995
$ws.Range("D9").Value = -1.7948619
$ws.Range("C10").Value = -112.25803
$ws.Range("D10").Value = -2.6328137
$ws.Range("D11").Value = 3.982335
$ws.Range("E11").Value = 22.365133
$ws.Range("C12").Value = -109.8891
$ws.Range("D12").Value = 0.099200339
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = -157.5
$ws.Range("D13").Value = 1.3451338
$ws.Range("E13").Value = 0.56036831
$ws.Range("C14").Value = -123.3891
$ws.Range("D14").Value = -1.076805
$ws.Range("C15").Value = -195.26745
$ws.Range("D15").Value = -0.63646691
$ws.Range("C16").Value = -258.87259
$ws.Range("D16").Value = -0.66413488
$ws.Range("C17").Value = -262.27099
$ws.Range("D17").Value = 1.2251036
$ws.Range("C18").Value = -325.87614
$ws.Range("D18").Value = 1.1974356
$ws.Range("C19").Value = -341.00632
$ws.Range("D19").Value = 0.092241832
$ws.Range("E19").Value = 0
$ws.Range("C20").Value = -116.53212
$ws.Range("D20").Value = 0.21360674
$ws.Range("C21").Value = -229.69731
$ws.Range("D21").Value = 0.11990981
$ws.Range("C22").Value = -355.05142
$ws.Range("D22").Value = 0.37351762
$ws.Range("C23").Value = -269.13237
$ws.Range("D23").Value = 0.25360781
$ws.Range("C24").Value = -310.53212
$ws.Range("D24").Value = 0.30006661
$ws.Range("C25").Value = -113.5781
$ws.Range("D25").Value = -0.10195156
$ws.Range("C26").Value = -161.51191
$ws.Range("D26").Value = 0.006460186
$ws.Range("C27").Value = -161.51191
$ws.Range("D27").Value = 0.006460186
$ws.Range("C28").Value = 274.60193
$ws.Range("D28").Value = 0.2981404
$ws.Range("C29").Value = -192.29401
$ws.Range("D29").Value = 0.057110836
$ws.Range("C30").Value = -176.81621
$ws.Range("D30").Value = -0.07148046399999999
$ws.Range("C31").Value = -68.68932100000001
$ws.Range("D31").Value = 0.02660047
$ws.Range("C32").Value = -123.60469
$ws.Range("D32").Value = 0.031136388
$ws.Range("C33").Value = -0.84466045
$ws.Range("D33").Value = 0.024700436
$ws.Range("C34").Value = -0.84466045
$ws.Range("D34").Value = 0.024700436
$ws.Range("C35").Value = -178.90811
$ws.Range("D35").Value = -0.062156925
$ws.Range("C36").Value = -178.90811
$ws.Range("D36").Value = -0.062156925
$ws.Range("C37").Value = -242.90811
$ws.Range("D37").Value = -0.034186309
$ws.Range("C38").Value = -242.90811
$ws.Range("D38").Value = -0.034186309
$ws.Range("C39").Value = -176.39531
$ws.Range("D39").Value = -0.020164518

$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = 104.18285
$ws.Range("B3").Value = 104.85521
$ws.Range("C3").Value = -0.34986073
$ws.Range("B4").Value = 101.62745
$ws.Range("C4").Value = 19.355457
$ws.Range("B5").Value = 102.66341
$ws.Range("C5").Value = 6.0649027
$ws.Range("B6").Value = 101.63979
$ws.Range("C6").Value = 3.5069326
$ws.Range("B7").Value = 117.38977
$ws.Range("C7").Value = 3.7781393
$ws.Range("B8").Value = 99.984578
$ws.Range("C8").Value = -13.676798
$ws.Range("B9").Value = 100.08378
$ws.Range("C9").Value = -6.9735627
$ws.Range("B10").Value = 100.86854
$ws.Range("C10").Value = 19.013937
$ws.Range("B11").Value = 99.006973
$ws.Range("C11").Value = 13.385639
$ws.Range("B12").Value = 100.23208
$ws.Range("C12").Value = 35.416403
$ws.Range("B13").Value = 100.20441
$ws.Range("C13").Value = 40.759235
$ws.Range("B14").Value = 100.32432
$ws.Range("C14").Value = 51.784706
$ws.Range("B15").Value = 100.44568
$ws.Range("C15").Value = 40.310752
$ws.Range("B16").Value = 100.8477
$ws.Range("C16").Value = 56.701319
$ws.Range("B17").Value = 100.74575
$ws.Range("C17").Value = 58.632147
$ws.Range("B18").Value = 100.80286
$ws.Range("C18").Value = 63.631791
$ws.Range("B19").Value = 100.82946
$ws.Range("C19").Value = 64.593442
$ws.Range("B20").Value = 100.67427
$ws.Range("C20").Value = 62.69892
$ws.Range("B21").Value = 100.61211
$ws.Range("C21").Value = 69.855244
$ws.Range("B22").Value = 100.85416
$ws.Range("C22").Value = 64.615403
$ws.Range("B23").Value = 100.834
$ws.Range("C23").Value = 76.61028399999999
$ws.Range("B24").Value = 100.57793
$ws.Range("C24").Value = 75.19922200000001
$ws.Range("B25").Value = 101.14584
$ws.Range("C25").Value = 42.422019
